$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Confusion-Matrix / Classification-Report / ROC-AUC / Accuracy / MSE / R2
# metrics (columns B:AE) produced by retraining + re-selecting each model (row 2:11).
$rows = @{}
$rows[2] = @(292, 89, 82, 286, 196, -7, -7, -210, 0.776566757493188, 0.7664041994750657, 0.7620320855614974, 0.7807486631016043, 0.7692307692307693, 0.7735099337748345, 374, 374, 0.7713903743315508, 748, 0.7714854784841269, 0.7713903743315509, 0.7713703515028019, 748, 0.771485478484127, 0.7713903743315508, 0.7713703515028019, 748, 0.8504461094111927, 0.7713903743315508, 0.2286096256684492, 0.08556149732620322)
$rows[3] = @(273, 117, 101, 257, 140, -16, -16, -172, 0.7178770949720671, 0.7, 0.6871657754010695, 0.7299465240641712, 0.7021857923497268, 0.7146596858638743, 374, 374, 0.7085561497326203, 748, 0.7089385474860335, 0.7085561497326203, 0.7084227391068005, 748, 0.7089385474860335, 0.7085561497326203, 0.7084227391068004, 748, 0.7893884583488233, 0.7085561497326203, 0.2914438502673797, -0.1657754010695187)
$rows[4] = @(282, 85, 92, 289, 204, 7, 7, -190, 0.7585301837270341, 0.7683923705722071, 0.7727272727272727, 0.7540106951871658, 0.7655629139072848, 0.7611336032388664, 374, 374, 0.7633689839572193, 748, 0.7634612771496205, 0.7633689839572193, 0.7633482585730755, 748, 0.7634612771496206, 0.7633689839572193, 0.7633482585730758, 748, 0.8420815579513283, 0.7633689839572193, 0.2366310160427808, 0.053475935828877)
$rows[5] = @(282, 85, 92, 289, 204, 7, 7, -190, 0.7585301837270341, 0.7683923705722071, 0.7727272727272727, 0.7540106951871658, 0.7655629139072848, 0.7611336032388664, 374, 374, 0.7633689839572193, 748, 0.7634612771496205, 0.7633689839572193, 0.7633482585730755, 748, 0.7634612771496206, 0.7633689839572193, 0.7633482585730758, 748, 0.8420815579513283, 0.7633689839572193, 0.2366310160427808, 0.053475935828877)
$rows[6] = @(328, 154, 46, 220, 66, -108, -108, -282, 0.8270676691729323, 0.6804979253112033, 0.5882352941176471, 0.8770053475935828, 0.6875, 0.7663551401869159, 374, 374, 0.732620320855615, 748, 0.7537827972420679, 0.732620320855615, 0.7269275700934579, 748, 0.7537827972420678, 0.732620320855615, 0.7269275700934579, 748, 0.8263104463953788, 0.732620320855615, 0.267379679144385, -0.06951871657754016)
$rows[7] = @(289, 100, 85, 274, 174, -15, -15, -204, 0.7632311977715878, 0.7429305912596401, 0.732620320855615, 0.7727272727272727, 0.747612551159618, 0.7575360419397117, 374, 374, 0.7526737967914439, 748, 0.753080894515614, 0.7526737967914439, 0.7525742965496649, 748, 0.7530808945156141, 0.7526737967914439, 0.7525742965496649, 748, 0.8337456032485917, 0.7526737967914439, 0.2473262032085561, 0.01069518716577544)
$rows[8] = @(289, 100, 85, 274, 174, -15, -15, -204, 0.7632311977715878, 0.7429305912596401, 0.732620320855615, 0.7727272727272727, 0.747612551159618, 0.7575360419397117, 374, 374, 0.7526737967914439, 748, 0.753080894515614, 0.7526737967914439, 0.7525742965496649, 748, 0.7530808945156141, 0.7526737967914439, 0.7525742965496649, 748, 0.8381495038462637, 0.7526737967914439, 0.2473262032085561, 0.01069518716577544)
$rows[9] = @(273, 104, 101, 270, 166, -3, -3, -172, 0.7277628032345014, 0.7241379310344828, 0.7219251336898396, 0.7299465240641712, 0.7248322147651006, 0.7270306258322237, 374, 374, 0.7259358288770054, 748, 0.725950367134492, 0.7259358288770054, 0.7259314202986622, 748, 0.725950367134492, 0.7259358288770054, 0.7259314202986621, 748, 0.8077297034516284, 0.7259358288770054, 0.2740641711229946, -0.09625668449197855)
$rows[10] = @(303, 105, 71, 269, 164, -34, -34, -232, 0.7911764705882353, 0.7426470588235294, 0.7192513368983957, 0.8101604278074866, 0.7535014005602241, 0.7749360613810742, 374, 374, 0.7647058823529411, 748, 0.7669117647058823, 0.7647058823529411, 0.7642187309706492, 748, 0.7669117647058823, 0.7647058823529411, 0.7642187309706492, 748, 0.8480296834338985, 0.7647058823529411, 0.2352941176470588, 0.05882352941176472)
$rows[11] = @(302, 108, 72, 266, 158, -36, -36, -230, 0.7869822485207101, 0.7365853658536585, 0.7112299465240641, 0.8074866310160428, 0.7471910112359551, 0.7704081632653061, 374, 374, 0.7593582887700535, 748, 0.7617838071871843, 0.7593582887700534, 0.7587995872506306, 748, 0.7617838071871843, 0.7593582887700535, 0.7587995872506307, 748, 0.8369913351825904, 0.7593582887700535, 0.2406417112299465, 0.03743315508021394)

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

Write-Host "Updated model evaluation metrics for rows 2-11 (B:AE)."
